# The document starts with an empty paragraph (centered, bold, size 32
# half-points i.e. 16pt, carried via the paragraph mark's rPr). The edit
# adds a single run containing just a space character to that paragraph,
# with explicit run formatting (bold + size 16pt) matching the paragraph
# mark's formatting.

$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)
$r = $p.Range
$r.Text = " "
$r.Font.Bold = $true
$r.Font.BoldBi = $true
$r.Font.Size = 16
$r.Font.SizeBi = 16
